$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "item" column (B) for every data row (2-544) with the
# placeholder "NA" value wherever it was previously blank or held a
# per-row unique string — the Hebrew_WG fix unifies these into a single
# shared "NA" string.
for ($i = 2; $i -le 544; $i++) {
    $ws.Range("B$i").Value = "NA"
}

# Move the view/selection to reflect where the editor ended up working.
[void]$ws.Range("B544").Select()
